# ---------------------------------------------------------------------------
# Edit script: reshape "combined Stats-this session" / "NL Stats-this
# session" from 5 players (rows 2-6) down to 3 players (rows 2-4), with new
# stat values, refresh the dependent charts' category/value ranges from
# $2:$*$6 to $2:$*$4, and flip the active tab back to the first sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# New data for rows 2-4 (identical on both worksheets). Row 5 and 6 are
# removed entirely (fewer players logged this session).
$playerRows = @(
    @{
        A = "Fish";  B = 20; C = 53.59; D = 33.59; E = 0
        F = 0.679;   G = 0;  H = 0;     I = 0.118;  J = 0.268; K = 0.19
        L = 0.19;    M = 0;  N = 0
        O = 112.17;  P = 66.58; Q = 168; R = 0.711
    },
    @{
        A = "Raymond"; B = 51; C = 24.39; D = -26.61; E = 0
        F = 0.466;     G = 0.31; H = 0.0057; I = 0.454; J = 0.126; K = 0.046
        L = 3.68;      M = 26; N = 56
        O = 41.7;      P = 58.54; Q = 174; R = 0.364
        T = "07/05/21"
    },
    @{
        A = "Scott"; B = 20; C = 13.02; D = -6.98; E = 0
        F = 0.527;   G = 0.24; H = 0.02; I = 0.513; J = 0.22; K = 0.08
        L = 3.16;    M = 16; N = 45
        O = 43.65;   P = 64.11; Q = 150; R = 0.364
    }
)

foreach ($ws in $wb.Worksheets) {

    # Drop the two trailing player rows (delete bottom-up so row numbers
    # of the rows still to be removed don't shift).
    $ws.Rows.Item(6).Delete()
    $ws.Rows.Item(5).Delete()

    # Write the refreshed per-player stats into rows 2-4. The "T" column
    # holds a free-text date label ("07/05/21") -- force text formatting
    # first so the COM layer doesn't auto-coerce the literal into a date
    # serial, then drop the formatting override again so the cell ends up
    # plain/unstyled (same as the source file).
    for ($r = 0; $r -lt $playerRows.Count; $r++) {
        $rowNum = $r + 2
        $rowData = $playerRows[$r]
        foreach ($col in $rowData.Keys) {
            $cell = $ws.Range("$col$rowNum")
            if ($col -eq "T") {
                $cell.NumberFormat = "@"
                $cell.Value2 = $rowData[$col]
                $cell.ClearFormats()
            } else {
                $cell.Value2 = $rowData[$col]
            }
        }
    }

    # Refresh every chart series on this sheet: the category axis and the
    # value series both referenced rows 2:6, now they should reference
    # rows 2:4.
    $chartCount = $ws.ChartObjects().Count
    for ($i = 1; $i -le $chartCount; $i++) {
        $chartObj = $ws.ChartObjects().Item($i)
        $chart = $chartObj.Chart
        $seriesCount = $chart.SeriesCollection().Count
        for ($j = 1; $j -le $seriesCount; $j++) {
            $ser = $chart.SeriesCollection().Item($j)
            $parts = $ser.Formula -split ','
            $catRef = $parts[1]
            $valRef = $parts[2]
            $newCat = $catRef.Replace("`$6", "`$4")
            $newVal = $valRef.Replace("`$6", "`$4")
            $ser.XValues = "=" + $newCat
            $ser.Values = "=" + $newVal
        }
    }
}

# The workbook previously opened with the second tab ("NL Stats-this
# session") active; switch back to the first tab.
$wb.Worksheets.Item(1).Activate()
